$d = $word.ActiveDocument

$replacements = @(
    @("83×39=3237", "62×59=3658"),
    @("13×79=1027", "84×40=3360"),
    @("47×95=4465", "90×89=8010"),
    @("83×56=4648", "54×74=3996"),
    @("65×15=975", "33×31=1023"),
    @("76×23=1748", "25×98=2450"),
    @("22×62=1364", "75×11=825"),
    @("61×86=5246", "81×62=5022"),
    @("92×33=3036", "58×70=4060"),
    @("23×47=1081", "11×23=253"),
    @("16×61=976", "42×71=2982"),
    @("91×87=7917", "49×27=1323"),
    @("15×31=465", "29×79=2291"),
    @("82×42=3444", "24×16=384"),
    @("93×12=1116", "87×94=8178"),
    @("23×34=782", "56×82=4592"),
    @("31×32=992", "54×70=3780"),
    @("73×73=5329", "85×99=8415"),
    @("44×70=3080", "97×75=7275"),
    @("38×74=2812", "76×13=988"),
    @("35×41=1435", "17×24=408"),
    @("19×53=1007", "14×94=1316"),
    @("86×51=4386", "60×45=2700"),
    @("43×63=2709", "91×74=6734"),
    @("76×54=4104", "37×34=1258")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
